# Update Vis Tool + Teams Data
# - "Games" sheet: append the completed game vs DAL (2023-12-01 / serial 45304)
#   that was previously listed as an upcoming game on the "Next" sheet.
# - "Next" sheet: remove that now-played game (row 2), shifting the
#   remaining upcoming games up by one row.

$wb = $excel.ActiveWorkbook

$wsGames = $wb.Worksheets.Item("Games")
$wsNext = $wb.Worksheets.Item("Next")

# New completed-game row (row 41) on the "Games" sheet.
$newRow = 41

$wsGames.Cells.Item($newRow, 1).Value = 40
$wsGames.Cells.Item($newRow, 2).Value = 45304
$wsGames.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"
$wsGames.Cells.Item($newRow, 3).Value = 1
$wsGames.Cells.Item($newRow, 4).Value = 118
$wsGames.Cells.Item($newRow, 5).Value = 89.7
$wsGames.Cells.Item($newRow, 6).Value = 0.57
$wsGames.Cells.Item($newRow, 7).Value = 9.5
$wsGames.Cells.Item($newRow, 8).Value = 35.9
$wsGames.Cells.Item($newRow, 9).Value = 0.354
$wsGames.Cells.Item($newRow, 10).Value = 131.6
$wsGames.Cells.Item($newRow, 11).Value = "DAL"
$wsGames.Cells.Item($newRow, 12).Value = 108
$wsGames.Cells.Item($newRow, 13).Value = 0.543
$wsGames.Cells.Item($newRow, 14).Value = 11.7
$wsGames.Cells.Item($newRow, 15).Value = 23.1
$wsGames.Cells.Item($newRow, 16).Value = 0.232
$wsGames.Cells.Item($newRow, 17).Value = 120.5
$wsGames.Cells.Item($newRow, 18).Value = 0
$wsGames.Cells.Item($newRow, 19).Value = 1

# Remove the now-played DAL game from the "Next" sheet (row 2); remaining
# rows shift up automatically.
$wsNext.Rows.Item(2).Delete()
